$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# River trend results update - May 2024
# Row 2 (Chlorophyll A)
$ws.Range("F2").Value = 0.663843034985196
$ws.Range("H2").Value = 0.672413793103448
$ws.Range("K2").Value = -2.65058055152395
$ws.Range("L2").Value = -14.2774411938331
$ws.Range("M2").Value = 8.189355470075339
$ws.Range("N2").Value = -1.68290828668187
$ws.Range("P2").Value = "As likely as not improving"

# Row 3 (Dissolved Oxygen Concentration)
$ws.Range("F3").Value = 0.203636393977046
$ws.Range("H3").Value = 0.862745098039216
$ws.Range("K3").Value = -0.07964953271027971
$ws.Range("L3").Value = -0.235012161645526
$ws.Range("M3").Value = 0.0845545742605892
$ws.Range("N3").Value = -0.708625735856581
$ws.Range("P3").Value = "Unlikely increasing"

# Row 4 (Ammoniacal Nitrogen / Total Nitrogen type row)
$ws.Range("E4").Value = "ok"
$ws.Range("F4").Value = 0.240440733425395
$ws.Range("J4").Value = 0.008999999999999999
$ws.Range("K4").Value = 0.0002996310465961
$ws.Range("L4").Value = -0.0003780024532736
$ws.Range("M4").Value = 0.0012811047389365
$ws.Range("N4").Value = 3.32923385106822

# Row 5 (Ammoniacal Nitrogen (NH4)) - now has too few non-censored values to analyse
$ws.Range("E5").Value = "< 5 Non-censored values"
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = 0.916666666666667
$ws.Range("H5").Value = 0.1875
$ws.Range("I5").Value = 5
$ws.Range("J5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("P5").Value = "Not Analysed improving"

# Row 6 (Nitrite Nitrogen (NO2))
$ws.Range("F6").Value = 0.998384516613228
$ws.Range("P6").Value = "Virtually certain improving"

# Row 7 (Nitrate Nitrogen (NO3))
$ws.Range("H7").Value = 0.0384615384615385
$ws.Range("I7").Value = 1

# Row 8 (pH)
$ws.Range("F8").Value = 0.032349197504054
$ws.Range("H8").Value = 0.788461538461538
$ws.Range("J8").Value = 7.865
$ws.Range("K8").Value = -0.0478243935367736
$ws.Range("L8").Value = -0.0854060824177271
$ws.Range("M8").Value = -0.0070090977074234
$ws.Range("N8").Value = -0.608066033525411
$ws.Range("P8").Value = "Extremely unlikely increasing"

# Row 9 (SIN (Soluble Inorganic nitrogen))
$ws.Range("F9").Value = 0.995257725331926
$ws.Range("G9").Value = 0.211538461538462
$ws.Range("H9").Value = 0.365384615384615
$ws.Range("I9").Value = 3
$ws.Range("J9").Value = 0.00675
$ws.Range("K9").Value = -0.0010870535714285
$ws.Range("L9").Value = -0.0014916328258204
$ws.Range("M9").Value = -0.0006202927021782
$ws.Range("N9").Value = -16.1044973544974
$ws.Range("P9").Value = "Virtually certain improving"

# Row 10
$ws.Range("F10").Value = 0.0002401380594218
$ws.Range("H10").Value = 0.529411764705882
$ws.Range("J10").Value = 140
$ws.Range("K10").Value = 8.790613718411549
$ws.Range("L10").Value = 5.16522093124125
$ws.Range("M10").Value = 12.7945238302441
$ws.Range("N10").Value = 6.27900979886539

# Row 11
$ws.Range("F11").Value = 0.418446339907522
$ws.Range("H11").Value = 0.88
$ws.Range("J11").Value = 11.24
$ws.Range("K11").Value = -0.013330291970803
$ws.Range("L11").Value = -0.102000617808653
$ws.Range("M11").Value = 0.0914970580767371
$ws.Range("N11").Value = -0.118596903654831
$ws.Range("P11").Value = "As likely as not increasing"

# Row 12
$ws.Range("F12").Value = 0.9978983410042001
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0.278350515463918
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = -0.0005482878963142
$ws.Range("L12").Value = -0.0008632355105724
$ws.Range("M12").Value = -0.0002468709211968
$ws.Range("N12").Value = -4.98443542103829
$ws.Range("P12").Value = "Virtually certain improving"

# Row 13
$ws.Range("F13").Value = 0.825180588708565
$ws.Range("G13").Value = 0.887640449438202
$ws.Range("H13").Value = 0.258426966292135
$ws.Range("I13").Value = 13
$ws.Range("P13").Value = "Likely improving"

# Row 14
$ws.Range("F14").Value = 0.90815807708855
$ws.Range("G14").Value = 0.752577319587629
$ws.Range("H14").Value = 0.0618556701030928
$ws.Range("P14").Value = "Very likely improving"

# Row 15 - now has too few unique values to analyse
$ws.Range("E15").Value = "< 3 unique values"
$ws.Range("F15").ClearContents()
$ws.Range("G15").Value = 0.958762886597938
$ws.Range("H15").Value = 0.0309278350515464
$ws.Range("J15").ClearContents()
$ws.Range("K15").ClearContents()
$ws.Range("L15").ClearContents()
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("P15").Value = "Not Analysed improving"

# Row 16
$ws.Range("F16").Value = 0.099424500971085
$ws.Range("H16").Value = 0.720430107526882
$ws.Range("J16").Value = 7.88
$ws.Range("K16").Value = -0.0162341999666984
$ws.Range("L16").Value = -0.0350003663974085
$ws.Range("M16").Value = 0.0035307275466673
$ws.Range("N16").Value = -0.206017766074853
$ws.Range("P16").Value = "Very unlikely increasing"

# Row 17
$ws.Range("E17").Value = "ok"
$ws.Range("F17").Value = 0.943025494138129
$ws.Range("G17").Value = 0.329896907216495
$ws.Range("H17").Value = 0.360824742268041
$ws.Range("K17").Value = -0.0001611871452484
$ws.Range("L17").Value = -0.0004784912855921
$ws.Range("N17").Value = -2.30267350354872
$ws.Range("P17").Value = "Very likely improving"
